$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'63.767.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'  +3.33%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'3.133.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  +1.97%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Formula = "'  -0.10%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'590.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  +1.87%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'146.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'  +3.65%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Formula = "'  -0.13%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Formula = "'3.125.55"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Formula = "'  +2.11%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Formula = "'  +1.75%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Formula = "'0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Formula = "'5.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "'  +3.60%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Formula = "'0.470"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "'  +0.61%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Formula = "'  +6.01%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Formula = "'36.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  +3.64%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Formula = "'  -0.60%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Formula = "'3.653.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'  +1.73%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Formula = "'  -0.96%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Formula = "'63.691.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'  +3.24%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Formula = "'3.129.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "'  +1.60%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Formula = "'466.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'  +4.03%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'14.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'  +3.97%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Formula = "'0.736"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "'  +1.01%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Formula = "'7.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'  +1.67%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Formula = "'13.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'  -3.43%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Formula = "'82.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "'  +0.81%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Formula = "'  +0.01%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Formula = "'9.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "'  +12.43%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Formula = "'  +2.39%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Formula = "'2.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Formula = "'  -0.87%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D31").Formula = "'6.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Formula = "'  +2.83%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Formula = "'27.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "'  +2.28%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Formula = "'  -2.59%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Formula = "'0.0₃0868"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Formula = "'  +9.21%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Formula = "'  +9.50%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Formula = "'  +2.01%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Formula = "'3.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Formula = "'  +13.34%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Formula = "'6.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Formula = "'  +1.72%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Formula = "'  +1.37%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Formula = "'450.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'  +5.45%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Formula = "'8.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "'  -0.54%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Formula = "'  +0.71%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Formula = "'2.904.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'  +3.94%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Formula = "'0.280"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "'  +4.15%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Formula = "'  +2.62%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Formula = "'  +6.35%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Formula = "'36.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'  +2.94%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Formula = "'125.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'  +0.84%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Formula = "'  +0.00%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Formula = "'  +0.50%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Formula = "'24.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "'  +3.91%  "
$ws.Range("E51").Style = "Normal"
